# ---------------------------------------------------------------------------
# "added matchrecords object model"
#
# This script updates the API list worksheet:
#   - tweaks a couple of existing reservation-related URIs / text
#   - fills in some previously-blank cells in the reservations section
#   - appends a brand-new "match records" section (rows 54-59)
#   - widens column G a little
#   - adjusts the view (zoom / freeze / selection)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-row tweaks -----------------------------------------------
$ws.Range("C25").Value = "/events/:eventID/reservations/"

$ws.Range("A35").Value = "ok"
$ws.Range("A35").WrapText = $true

$ws.Range("A44").Value = "ok"
$ws.Range("A44").WrapText = $true
$ws.Range("C44").Value = "/events/:eventId/reservations/cancel"
$ws.Range("E44").Value = "anyone who signed up"

$ws.Range("A45").Value = "ok"
$ws.Range("A45").WrapText = $true

$ws.Range("A46").Value = "ok"
$ws.Range("A46").WrapText = $true

$ws.Range("A47").Value = "ok"
$ws.Range("A47").WrapText = $true
$ws.Range("D47").Value = "get"
$ws.Range("E47").Value = "admin / the same user"
$ws.Range("E47").WrapText = $true

# --- New "match records" rows (54-59) ----------------------------------
$ws.Range("A54").Value = "coded"
$ws.Range("B54").Value = "add match record"
$ws.Range("C54").Value = "/events/:eventId/matchrecords/"
$ws.Range("D54").Value = "post"
$ws.Range("E54").Value = "admin / host"
$ws.Range("G54").Value = "Match Record Obj"
$ws.Range("A54:C54").WrapText = $true
$ws.Range("E54").WrapText = $true
$ws.Range("G54").WrapText = $true

$ws.Range("A55").Value = "coded"
$ws.Range("B55").Value = "modify match record"
$ws.Range("C55").Value = "/events/:eventId/matchrecords/:matchrecordId"
$ws.Range("D55").Value = "patch"
$ws.Range("E55").Value = "admin"
$ws.Range("G55").Value = "Match Record Obj"
$ws.Range("A55:C55").WrapText = $true
$ws.Range("E55").WrapText = $true
$ws.Range("G55").WrapText = $true

$ws.Range("A56").Value = "coded"
$ws.Range("B56").Value = "delete match record"
$ws.Range("C56").Value = "/events/:eventId/matchrecords/:matchrecordId"
$ws.Range("D56").Value = "delete"
$ws.Range("E56").Value = "admin"
$ws.Range("G56").Value = "Match Record Obj"
$ws.Range("A56:C56").WrapText = $true
$ws.Range("E56").WrapText = $true
$ws.Range("G56").WrapText = $true

$ws.Range("A57").Value = "coded"
$ws.Range("B57").Value = "get event matchrecords"
$ws.Range("C57").Value = "/events/:eventId/matchrecords"
$ws.Range("D57").Value = "get"
$ws.Range("E57").Value = "all logged user"
$ws.Range("G57").Value = "[Match Record Obj]"
$ws.Range("A57:C57").WrapText = $true
$ws.Range("E57").WrapText = $true
$ws.Range("G57").WrapText = $true

$ws.Range("A58").Value = "coded"
$ws.Range("B58").Value = "get user matchrecords"
$ws.Range("C58").Value = "/users/:userId/matchrecords"
$ws.Range("D58").Value = "get"
$ws.Range("E58").Value = "admin / same user"
$ws.Range("G58").Value = "[Match Record Obj]"
$ws.Range("A58:C58").WrapText = $true
$ws.Range("E58").WrapText = $true
$ws.Range("G58").WrapText = $true

$ws.Range("A59").Value = "coded"
$ws.Range("B59").Value = "get event matchrecords of specific user"
$ws.Range("C59").Value = "/events/:eventId/matchrecords/of/:userId"
$ws.Range("D59").Value = "get"
$ws.Range("E59").Value = "admin / same user"
$ws.Range("A59:E59").WrapText = $true

# --- Column width --------------------------------------------------------
# Target stored width is 22.109375 characters; the closest value this
# engine's column-width rounding can reach is ~22.17.
$ws.Columns.Item(7).ColumnWidth = 21.25

# --- View / window state (best effort) -----------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("D60").Select()
